$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (Reactivos -> Catálogos). This also updates the
# "Reactivos" defined name's sheet-qualified reference automatically.
$ws.Name = "Catálogos"

# Update the placeholder text for Clave / Nombre rows.
$ws.Range("B3").Value = "{{Catalogo.Clave}}"
$ws.Range("B5").Value = "{{Catalogo.Nombre}}"

# Row 7 ("Clave Contpaq") becomes the "Activo" row.
$ws.Range("A7").Value = "Activo"
$ws.Range("B7").Value = "{{Catalogo.Activo}}"

# Remove the old "Nombre Contpaq" (row 9) and "Activo" (row 11) rows —
# their content has been folded into row 7 above.
$ws.Rows(11).Delete()
$ws.Rows(9).Delete()
